$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D"  = 0.101
    "G"  = 0.1900797373358349
    "H"  = 0.1900797373358349
    "I"  = 0.1304500681763358
    "J"  = 0.06522503408816788
    "K"  = -21.2
    "L"  = -0.0124296435272045
    "M"  = 64
    "N"  = 0.08695652173913043
    "O"  = -3.018867924528302
    "P"  = 32.5
    "Q"  = 0.04415760869565218
    "R"  = -1.533018867924528
    "S"  = 31.5
    "T"  = 0.4921875
    "U"  = 384.2
    "V"  = 0.5220108695652174
    "W"  = -0.03193733052124134
    "X"  = 0.3485015326537834
    "Y"  = -0.3804388631750247
    "Z"  = 0.942973324109291
    "AA" = 0.06150546720926149
    "AB" = 0.1686111665860797
    "AC" = -0.1071056993768182
    "AD" = 1506.5
    "AE" = 0.446818592208664
    "AF" = 1506.946818592209
    "AG" = 1122.746818592209
    "AH" = 0.6718602537076861
    "AI" = 0.4816454976917132
    "AJ" = 0.6040342920089369
    "AK" = 0.4090827713291158
    "AL" = 44.5
    "AM" = 44.5
    "AN" = 6.129340683931077
    "AO" = 4.991011235955056
    "AP" = 4.56800381875301
    "AQ" = 4.991011235955056
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range("${col}2").Value = $val
    $ws.Range("${col}3").Value = $val
}
